$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowFV($row) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += , ($ws.Range($c + $row).Value2)
    }
    return $vals
}

function Set-RowFV($row, $vals) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowFV $rowA
    $b = Get-RowFV $rowB
    Set-RowFV $rowA $b
    Set-RowFV $rowB $a
}

# Several match rows were re-ordered (same kickoff date/time slot, home/away
# pair flipped in the source listing) -- swap the F:V (match detail) content
# between each pair while leaving the Indice/pais/torneio/temporada/data
# columns (A:E) untouched.
Swap-Rows 22 23
Swap-Rows 76 77
Swap-Rows 96 97
Swap-Rows 107 108
Swap-Rows 139 140

# Append the two newly scraped matches as rows 160 and 161.
$ws.Range("A159:V159").Copy()
$ws.Range("A160:V161").PasteSpecial(-4122)

$newRows = @(
    @{
        A = 159; B = "indonesia"; C = "liga-1"; D = "2023-2024"; E = 45234.375
        F = "Persikabo 1973"; G = 2; H = "RANS Nusantara"; I = 1
        J = 2.77; K = "02/11/2023 21:12"
        L = 2.43; M = "04/11/2023 08:53"
        N = 3.23; O = "02/11/2023 21:12"
        P = 3.43; Q = "04/11/2023 08:52"
        R = 2.38; S = "02/11/2023 21:12"
        T = 2.82; U = "04/11/2023 08:53"
        V = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-rans-nusantara/UezuT8sm/"
    },
    @{
        A = 160; B = "indonesia"; C = "liga-1"; D = "2023-2024"; E = 45234.54166666666
        F = "Persita"; G = 2; H = "Barito Putera"; I = 2
        J = 2.6; K = "03/11/2023 01:12"
        L = 2.49; M = "04/11/2023 12:53"
        N = 2.95; O = "03/11/2023 01:12"
        P = 3.21; Q = "04/11/2023 12:53"
        R = 2.71; S = "03/11/2023 01:12"
        T = 2.91; U = "04/11/2023 12:53"
        V = "https://www.betexplorer.com/football/indonesia/liga-1/persita-ps-barito-putera/OIZrSSdg/"
    }
)

$targetRow = 160
foreach ($rowData in $newRows) {
    $ws.Range("A$targetRow").Value = $rowData.A
    $ws.Range("B$targetRow").Value = $rowData.B
    $ws.Range("C$targetRow").Value = $rowData.C
    $ws.Range("D$targetRow").Value = $rowData.D
    $ws.Range("E$targetRow").Value = $rowData.E
    $ws.Range("F$targetRow").Value = $rowData.F
    $ws.Range("G$targetRow").Value = $rowData.G
    $ws.Range("H$targetRow").Value = $rowData.H
    $ws.Range("I$targetRow").Value = $rowData.I
    $ws.Range("J$targetRow").Value = $rowData.J
    $ws.Range("K$targetRow").Value = $rowData.K
    $ws.Range("L$targetRow").Value = $rowData.L
    $ws.Range("M$targetRow").Value = $rowData.M
    $ws.Range("N$targetRow").Value = $rowData.N
    $ws.Range("O$targetRow").Value = $rowData.O
    $ws.Range("P$targetRow").Value = $rowData.P
    $ws.Range("Q$targetRow").Value = $rowData.Q
    $ws.Range("R$targetRow").Value = $rowData.R
    $ws.Range("S$targetRow").Value = $rowData.S
    $ws.Range("T$targetRow").Value = $rowData.T
    $ws.Range("U$targetRow").Value = $rowData.U
    $ws.Range("V$targetRow").Value = $rowData.V
    $targetRow++
}

Write-Host "Done. Dimension:" $ws.UsedRange.Rows.Count
